$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Add row 13 to the Logs sheet
$ws1.Cells.Item(13, 1).Value = "Klacht over levering"
$ws1.Cells.Item(13, 2).Value = "mailmind.test@zohomail.eu"
$ws1.Cells.Item(13, 3).Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$ws1.Cells.Item(13, 4).Value = "Klacht"
$ws1.Cells.Item(13, 6).Value = "2025-06-18 14:00:12"
$ws1.Cells.Item(13, 7).Value = "Nee"

# Add row 14 to the Logs sheet
$ws1.Cells.Item(14, 1).Value = "Vragen over samenwerking"
$ws1.Cells.Item(14, 2).Value = "mailmind.test@zohomail.eu"
$ws1.Cells.Item(14, 3).Value = "Kunnen we samenwerken aan een nieuw project?"
$ws1.Cells.Item(14, 4).Value = "Overig"
$ws1.Cells.Item(14, 6).Value = "2025-06-18 14:00:13"
$ws1.Cells.Item(14, 7).Value = "Nee"

# Extend conditional formatting ranges to cover the new rows (D2:D12 -> D2:D14, G2:G12 -> G2:G14)
$dFormats = $ws1.Range("D2:D12").FormatConditions
for ($i = 1; $i -le $dFormats.Count; $i++) {
    $dFormats.Item($i).ModifyAppliesToRange($ws1.Range("D2:D14"))
}

$gFormats = $ws1.Range("G2:G12").FormatConditions
for ($i = 1; $i -le $gFormats.Count; $i++) {
    $gFormats.Item($i).ModifyAppliesToRange($ws1.Range("G2:G14"))
}

# Update the Dashboard counts: Overig 4 -> 5, Klacht 3 -> 4
$ws2.Range("B2").Value = 5
$ws2.Range("B3").Value = 4
